$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.134.76'
$ws.Range("E2").Value = '  -0.69%  '

$ws.Range("D3").Value = '3.157.88'
$ws.Range("E3").Value = '  -0.04%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.28'
$ws.Range("E5").Value = '  -1.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.92'
$ws.Range("E6").Value = '  -2.58%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = '3.150.79'
$ws.Range("E8").Value = '  -0.14%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.525'
$ws.Range("E9").Value = '  -0.78%  '

$ws.Range("E10").Value = '  -1.35%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.39'
$ws.Range("E11").Value = '  -1.44%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.467'
$ws.Range("E12").Value = '  -1.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000255'
$ws.Range("E13").Value = '  -1.35%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.09'
$ws.Range("E14").Value = '  -1.86%  '

$ws.Range("D15").Value = '3.669.55'
$ws.Range("E15").Value = '  -0.25%  '

$ws.Range("E16").Value = '  +2.45%  '

$ws.Range("D17").Value = '64.123.83'
$ws.Range("E17").Value = '  -0.76%  '

$ws.Range("D18").Value = '3.143.54'
$ws.Range("E18").Value = '  -0.59%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.86'
$ws.Range("E19").Value = '  -1.01%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '490.36'
$ws.Range("E20").Value = '  +2.46%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.70'
$ws.Range("E21").Value = '  -0.20%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.713'
$ws.Range("E22").Value = '  -0.90%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.66'
$ws.Range("E23").Value = '  -4.50%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '87.83'
$ws.Range("E24").Value = '  +4.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.33'
$ws.Range("E25").Value = '  -3.08%  '

$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("E27").Value = '  -2.34%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.25'
$ws.Range("E28").Value = '  -3.61%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.05'
$ws.Range("E29").Value = '  +0.96%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.07'
$ws.Range("E30").Value = '  -0.66%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '27.71'
$ws.Range("E31").Value = '  +4.53%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.112'
$ws.Range("E32").Value = '  -5.35%  '

$ws.Range("E33").Value = '  -0.07%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.66'
$ws.Range("E34").Value = '  -1.87%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.11'
$ws.Range("E35").Value = '  -2.48%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.04'
$ws.Range("E36").Value = '  +0.70%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.76'
$ws.Range("E37").Value = '  -0.39%  '

$ws.Range("D38").Value = '0.0₃0749'
$ws.Range("E38").Value = '  -3.76%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.99'
$ws.Range("E39").Value = '  -6.80%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '436.50'
$ws.Range("E40").Value = '  -5.21%  '

$ws.Range("E41").Value = '  -0.47%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.120'
$ws.Range("E42").Value = '  -1.35%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.32'
$ws.Range("E43").Value = '  -1.10%  '

$ws.Range("D44").Value = '2.931.10'
$ws.Range("E44").Value = '  +2.99%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.260'
$ws.Range("E45").Value = '  -3.37%  '

$ws.Range("E46").Value = '  -5.42%  '

$ws.Range("E47").Value = '  -2.64%  '

$ws.Range("E48").Value = '  -0.06%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '26.01'
$ws.Range("E49").Value = '  -2.20%  '

$ws.Range("E50").Value = '  -0.07%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.71'
$ws.Range("E51").Value = '  +0.32%  '
